$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data table contents (rows 16-29), grouped by worker, periods descending.
# Columns: C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora, F = Valor Mora
$data = @(
    @{ Row = 16; Doc = "9101123"; Nombre = "DAVID DE JESUS MARTINEZ URZOLA"; Periodo = "2406"; Valor = 16000 },
    @{ Row = 17; Doc = "9101123"; Nombre = "DAVID DE JESUS MARTINEZ URZOLA"; Periodo = "2405"; Valor = 40000 },
    @{ Row = 18; Doc = "9101123"; Nombre = "DAVID DE JESUS MARTINEZ URZOLA"; Periodo = "2404"; Valor = 40000 },
    @{ Row = 19; Doc = "9101123"; Nombre = "DAVID DE JESUS MARTINEZ URZOLA"; Periodo = "2403"; Valor = 40000 },
    @{ Row = 20; Doc = "9101123"; Nombre = "DAVID DE JESUS MARTINEZ URZOLA"; Periodo = "2402"; Valor = 40000 },
    @{ Row = 21; Doc = "9101123"; Nombre = "DAVID DE JESUS MARTINEZ URZOLA"; Periodo = "2401"; Valor = 46400 },
    @{ Row = 22; Doc = "9101123"; Nombre = "DAVID DE JESUS MARTINEZ URZOLA"; Periodo = "2312"; Valor = 46400 },
    @{ Row = 23; Doc = "9295903"; Nombre = "JOSE DEL CARMEN OSPINO CARO";    Periodo = "2406"; Valor = 18560 },
    @{ Row = 24; Doc = "9295903"; Nombre = "JOSE DEL CARMEN OSPINO CARO";    Periodo = "2405"; Valor = 46400 },
    @{ Row = 25; Doc = "9295903"; Nombre = "JOSE DEL CARMEN OSPINO CARO";    Periodo = "2404"; Valor = 46400 },
    @{ Row = 26; Doc = "9295903"; Nombre = "JOSE DEL CARMEN OSPINO CARO";    Periodo = "2403"; Valor = 46400 },
    @{ Row = 27; Doc = "9295903"; Nombre = "JOSE DEL CARMEN OSPINO CARO";    Periodo = "2402"; Valor = 46400 },
    @{ Row = 28; Doc = "9295903"; Nombre = "JOSE DEL CARMEN OSPINO CARO";    Periodo = "2401"; Valor = 46400 },
    @{ Row = 29; Doc = "9295903"; Nombre = "JOSE DEL CARMEN OSPINO CARO";    Periodo = "2312"; Valor = 46400 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = $entry.Doc
    $ws.Cells.Item($r, 4).Value = $entry.Nombre
    $ws.Cells.Item($r, 5).Value = $entry.Periodo
    $ws.Cells.Item($r, 6).Value = $entry.Valor
}
